$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "nan" value out of C9 (the "polite_expressions" cell for
# the existing "Thank the reviewer..." annotation) - it becomes blank.
$ws.Cells.Item(9, 3).Value = ""

# Append a new annotation row (row 10) for parisk.
$ws.Cells.Item(10, 1).Value = "parisk"
$ws.Cells.Item(10, 2).Value = 2
$ws.Cells.Item(10, 3).Value = "nan"
$ws.Cells.Item(10, 4).Value = "DFT"
$ws.Cells.Item(10, 5).Value = "THE"
$ws.Cells.Item(10, 6).Value = "a5228610-fe6d-4383-b598-a7c34c3b8714"
$ws.Cells.Item(10, 7).Value = "HyRnez-RW_annotated.xlsx"
$ws.Cells.Item(10, 8).Value = "Why is this result not compared to in Table 1?"
